$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new "Dinh danh SPHH" row of data right after the existing
#     table rows (was A4:F32, now one more data row at 33). It will find its
#     correct alphabetical slot once the table body is (re)sorted below,
#     exactly like the diff shows a new row for the identifier-scanner
#     feature added to the icon/control inventory. ---
$ws.Range("A33").Value = "VnsErp2025"
$ws.Range("B33").Value = "FormMain"
$ws.Range("C33").Value = "XuatKhoRibbonPageGroup"
$ws.Range("D33").Value = "DinhDanhSpHhBarButtonItem"
$ws.Range("E33").Value = "Định danh SPHH"
$ws.Range("F33").Value = "brand-image.svg"

# --- Re-sort the table body (A4:F33) ascending by column C (MODULE group),
#     same as re-applying Table1's sort after adding the new row. This is
#     what pushes "XuatNhapTonKhoRibbonPage" (was row 18) down to the very
#     bottom (row 33), since it now sorts last alphabetically. ---
$table = $ws.ListObjects.Item("Table1")
$table.DataBodyRange.Sort($ws.Range("C4:C33"), 1)

# --- Restore the selection to what was left active on save. ---
$ws.Range("C39").Select() | Out-Null
